$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.758.60"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "1.625.01"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").Value = "'19.40"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "1.853.58"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13").Value = "1.637.79"
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("D14").Value = "'4.07"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").Value = "'65.10"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "26.738.87"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "'232.95"
$ws.Range("E18").Value = "  +9.58%  "
$ws.Range("D19").Value = "'7.77"
$ws.Range("E19").Value = "  +5.17%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'4.41"
$ws.Range("E22").Value = "  +3.23%  "
$ws.Range("D23").Value = "'2.26"
$ws.Range("E23").Value = "  +3.58%  "
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").Value = "'145.89"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'7.08"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("D29").Value = "'15.68"
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").Value = "'3.26"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("D33").Value = "1.471.90"
$ws.Range("E33").Value = "  +9.85%  "
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("E39").Value = "  +2.41%  "
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "'2.21"
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("D43").Value = "'0.957"
$ws.Range("E43").Value = "  -2.50%  "
$ws.Range("D44").Value = "1.764.42"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("D45").Value = "'0.767"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "'62.19"
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("D47").Value = "'88.46"
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("D48").Value = "'1.51"
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("D50").Value = "'0.0966"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").Value = "'7.47"
$ws.Range("E51").Value = "  +1.84%  "
